$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Removed stop words" note from H4 (it moves down to H16,
# with updated text, added further below).
$ws.Range("H4").ClearContents()

# --- New workflow-status notes in column F / legend in column K ---
# Write order controls the shared-string table build order, so we add the
# brand-new strings in the same sequence they first appear in the target
# workbook (idx 9..16).

# idx 9 / 10: the (reworded) stray notes that move to H16:H18
$ws.Range("H16").Value = "Removed podcast stop words"
$ws.Range("H17").Value = "removed short words"

# idx 11: "1. Preprocessed"
$ws.Range("F4").Value = "1. Preprocessed"
# idx 12: "2. Model trained"
$ws.Range("K5").Value = "2. Model trained"
# idx 13: "3. Categories Mapped to iTunes"
$ws.Range("F3").Value = "3. Categories Mapped to iTunes"
# idx 14: "4. Log-Likelihood calculated"
$ws.Range("K7").Value = "4. Log-Likelihood calculated"
# idx 15: long explanatory note
$ws.Range("H18").Value = "iTunes categories weren't well set up. Usubs found junk and then better categorized existing categories. For example …"
# idx 16: "(preprocessing)"
$ws.Range("G8").Value = "(preprocessing)"

# Remaining cells reusing already-created strings
$ws.Range("K4").Value = "1. Preprocessed"
$ws.Range("K6").Value = "3. Categories Mapped to iTunes"
$ws.Range("F6").Value = "2. Model trained"
$ws.Range("F7").Value = "1. Preprocessed"

# Selection moves to G4
$ws.Range("G4").Select()

# Data validation dropdown on F3 and F5:F8 (F4 already has a fixed value and
# is intentionally excluded), sourced from the legend list in column K.
$fullRange = $ws.Range("F3:F8")
$fullRange.Validation.Add(3, 1, 1, "=`$K`$4:`$K`$10")
$ws.Range("F4").Validation.Delete()
